$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.735.93"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.634.34"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.258"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0636"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.860.30"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.636.39"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "25.769.63"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "1.123.74"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.802"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "1.769.61"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "
